$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Егор Барсуков): several 1s become 0 — subsequent parsing
# apparently detected no qualifying answers for this entrant.
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0

# Row 3: was a placeholder/blank entrant (Фёдор Самохин with all
# zeros) — relabel it as the editable-variant proto entrant "Сергей
# Цыкура" (its score columns were already all zero, so they're left
# untouched).
$ws.Range("A3").Value = "Сергей"
$ws.Range("B3").Value = "Цыкура"

# Row 4: relabel Егор Барсуков -> Фёдор Самохин and zero out the
# scoring columns that used to hold this entrant's 1s.
$ws.Range("A4").Value = "Фёдор"
$ws.Range("B4").Value = "Самохин"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("AD4").Value = 0
$ws.Range("AE4").Value = 0

# Rows 5 and 6 (the old Сергей Цыкура / Фёдор Самохин rows) are gone
# now that their data was folded into rows 3 and 4 above — delete the
# now-duplicate trailing rows (delete row 6 first so row numbers stay
# stable for the second delete).
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
